# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / Leve profit figures across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row98
$ws.Range("H98").Value = 641.4054
$ws.Range("I98").Value = 501.07693
$ws.Range("J98").Value = 973.0909
$ws.Range("K98").Value = 501.07693
$ws.Range("L98").Value = 973.0909
$ws.Range("M98").Value = 996.9230700000001
$ws.Range("N98").Value = -3969.0909

# ALC!row116
$ws.Range("H116").Value = 78820.36
$ws.Range("I116").Value = 84575.766
$ws.Range("K116").Value = 84575.766
$ws.Range("M116").Value = -81133.766

# ALC!row122
$ws.Range("H122").Value = 641.4054
$ws.Range("I122").Value = 501.07693
$ws.Range("J122").Value = 973.0909
$ws.Range("K122").Value = 1503.23079
$ws.Range("L122").Value = 2919.2727
$ws.Range("M122").Value = 946.7692099999999
$ws.Range("N122").Value = -7819.2727

# ALC!row132
$ws.Range("H132").Value = 3067.25
$ws.Range("I132").Value = 1787.2433
$ws.Range("J132").Value = 7372.727
$ws.Range("K132").Value = 5361.7299
$ws.Range("L132").Value = 22118.181
$ws.Range("M132").Value = -2831.7299
$ws.Range("N132").Value = -27178.181

# ALC!row141
$ws.Range("H141").Value = 2150
$ws.Range("I141").Value = 2150
$ws.Range("K141").Value = 6450
$ws.Range("M141").Value = -1270

$ws = $wb.Worksheets.Item("ARM")
# ARM!row74
$ws.Range("H74").Value = 1958.909
$ws.Range("I74").Value = 1027.0968
$ws.Range("J74").Value = 4180.923
$ws.Range("K74").Value = 1027.0968
$ws.Range("L74").Value = 4180.923
$ws.Range("M74").Value = -153.0968
$ws.Range("N74").Value = -5928.923

# ARM!row77
$ws.Range("H77").Value = 1958.909
$ws.Range("I77").Value = 1027.0968
$ws.Range("J77").Value = 4180.923
$ws.Range("K77").Value = 5135.484
$ws.Range("L77").Value = 20904.615
$ws.Range("M77").Value = -767.4840000000004
$ws.Range("N77").Value = -29640.615

# ARM!row139
$ws.Range("H139").Value = 40151.668
$ws.Range("J139").Value = 40151.668
$ws.Range("L139").Value = 40151.668
$ws.Range("N139").Value = -50431.668

$ws = $wb.Worksheets.Item("BSM")
# BSM!row81
$ws.Range("H81").Value = 22780
$ws.Range("J81").Value = 22780
$ws.Range("L81").Value = 22780
$ws.Range("N81").Value = -24902

# BSM!row84
$ws.Range("H84").Value = 22780
$ws.Range("J84").Value = 22780
$ws.Range("L84").Value = 68340
$ws.Range("N84").Value = -78948

# BSM!row135
$ws.Range("H135").Value = 26156
$ws.Range("J135").Value = 26926.666
$ws.Range("L135").Value = 26926.666
$ws.Range("N135").Value = -37066.666

$ws = $wb.Worksheets.Item("CRP")
# CRP!row31
$ws.Range("H31").Value = 2512.3728
$ws.Range("I31").Value = 1488.5555
$ws.Range("J31").Value = 4114.8696
$ws.Range("K31").Value = 1488.5555
$ws.Range("L31").Value = 4114.8696
$ws.Range("M31").Value = -1193.5555
$ws.Range("N31").Value = -4704.8696

# CRP!row34
$ws.Range("H34").Value = 2512.3728
$ws.Range("I34").Value = 1488.5555
$ws.Range("J34").Value = 4114.8696
$ws.Range("K34").Value = 1488.5555
$ws.Range("L34").Value = 4114.8696
$ws.Range("M34").Value = -1286.5555
$ws.Range("N34").Value = -4518.8696

# CRP!row99
$ws.Range("H99").Value = 65231.125
$ws.Range("I99").Value = 102447.4
$ws.Range("J99").Value = 3204
$ws.Range("K99").Value = 102447.4
$ws.Range("L99").Value = 3204
$ws.Range("M99").Value = -100949.4
$ws.Range("N99").Value = -6200

# CRP!row126
$ws.Range("H126").Value = 65231.125
$ws.Range("I126").Value = 102447.4
$ws.Range("J126").Value = 3204
$ws.Range("K126").Value = 307342.2
$ws.Range("L126").Value = 9612
$ws.Range("M126").Value = -304872.2
$ws.Range("N126").Value = -14552

# CRP!row132
$ws.Range("H132").Value = 2517.6
$ws.Range("I132").Value = 1386.1818
$ws.Range("J132").Value = 3406.5715
$ws.Range("K132").Value = 4158.5454
$ws.Range("L132").Value = 10219.7145
$ws.Range("M132").Value = -1628.5454
$ws.Range("N132").Value = -15279.7145

# CRP!row134
$ws.Range("H134").Value = 1904.675
$ws.Range("I134").Value = 1177.0526
$ws.Range("J134").Value = 2563
$ws.Range("K134").Value = 3531.1578
$ws.Range("L134").Value = 7689
$ws.Range("M134").Value = -996.1578
$ws.Range("N134").Value = -12759

$ws = $wb.Worksheets.Item("CUL")
# CUL!row68
$ws.Range("H68").Value = 437.84616
$ws.Range("I68").Value = 476.8889
$ws.Range("J68").Value = 350
$ws.Range("K68").Value = 1430.6667
$ws.Range("L68").Value = 1050
$ws.Range("M68").Value = -619.6667
$ws.Range("N68").Value = -2672

# CUL!row71
$ws.Range("H71").Value = 437.84616
$ws.Range("I71").Value = 476.8889
$ws.Range("J71").Value = 350
$ws.Range("K71").Value = 4292.0001
$ws.Range("L71").Value = 3150
$ws.Range("M71").Value = -236.0001000000002
$ws.Range("N71").Value = -11262

# CUL!row75
$ws.Range("H75").Value = 1696.5333
$ws.Range("I75").Value = 582.6667
$ws.Range("J75").Value = 1975
$ws.Range("K75").Value = 1748.0001
$ws.Range("L75").Value = 5925
$ws.Range("M75").Value = -750.0001
$ws.Range("N75").Value = -7921

# CUL!row78
$ws.Range("H78").Value = 1696.5333
$ws.Range("I78").Value = 582.6667
$ws.Range("J78").Value = 1975
$ws.Range("K78").Value = 5244.0003
$ws.Range("L78").Value = 17775
$ws.Range("M78").Value = -252.0002999999997
$ws.Range("N78").Value = -27759

# CUL!row86
$ws.Range("H86").Value = 150
$ws.Range("I86").Value = 150
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 450
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 736
$ws.Range("N86").ClearContents()

# CUL!row89
$ws.Range("H89").Value = 150
$ws.Range("I89").Value = 150
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1350
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 4578
$ws.Range("N89").ClearContents()

# CUL!row114
$ws.Range("H114").Value = 1003.6818
$ws.Range("I114").Value = 482.4
$ws.Range("J114").Value = 1157
$ws.Range("K114").Value = 1447.2
$ws.Range("L114").Value = 3471
$ws.Range("M114").Value = 1806.8
$ws.Range("N114").Value = -9979

# CUL!row126
$ws.Range("H126").Value = 10926.667
$ws.Range("I126").Value = 12905.714
$ws.Range("K126").Value = 38717.142
$ws.Range("M126").Value = -33777.142

$ws = $wb.Worksheets.Item("GSM")
# GSM!row126
$ws.Range("H126").Value = 2292
$ws.Range("I126").Value = 1825.2307
$ws.Range("J126").Value = 2898.8
$ws.Range("K126").Value = 5475.6921
$ws.Range("L126").Value = 8696.400000000001
$ws.Range("M126").Value = -3005.6921
$ws.Range("N126").Value = -13636.4

$ws = $wb.Worksheets.Item("LTW")
# LTW!row40
$ws.Range("H40").Value = 3020
$ws.Range("I40").Value = 2823.6365
$ws.Range("J40").Value = 4100
$ws.Range("K40").Value = 2823.6365
$ws.Range("L40").Value = 4100
$ws.Range("M40").Value = -2687.6365
$ws.Range("N40").Value = -4372

# LTW!row121
$ws.Range("H121").Value = 38710
$ws.Range("J121").Value = 38710
$ws.Range("L121").Value = 38710
$ws.Range("N121").Value = -42204

